$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New CF rows for specific CO2 flows (to account for NETs)
$rows = @(
    @{A="Carbon dioxide, in air"; B="natural resource::in air"; C="-2.076923076923077E-13"},
    @{A="Carbon dioxide, non-fossil, resource correction"; B="natural resource::in air"; C="-2.076923076923077E-13"},
    @{A="Carbon dioxide, non-fossil"; B="air::lower stratosphere + upper troposphere"; C="2.076923076923077E-13"},
    @{A="Carbon dioxide, non-fossil"; B="air::non-urban air or from high stacks"; C="2.076923076923077E-13"},
    @{A="Carbon dioxide, non-fossil"; B="air::unspecified"; C="2.076923076923077E-13"},
    @{A="Carbon dioxide, non-fossil"; B="air::urban air close to ground"; C="2.076923076923077E-13"},
    @{A="Carbon dioxide, non-fossil"; B="air::non-urban air or from high stacks"; C="2.076923076923077E-13"}
)

$startRow = 78
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i].A
    $ws.Cells.Item($r, 2).Value = $rows[$i].B
    $ws.Cells.Item($r, 3).Value = [double]$rows[$i].C
}

# Restore the AutoFilter defined name (hidden, sheet-scoped) over the new data extent
$name = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$C`$84")
$name.Visible = $false

# Update the selected cell as recorded in the saved view
$ws.Range("B5").Select()
